# Fill in the "Average loss" column (F) for the per-configuration data rows
# (rows 2-64), matching the format already used by existing F-column cells
# (fillId=5 / borderId=1 "highlight" style), and update the active
# selection to F1 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing F-column number formatting (fill + border) onto the
# whole F2:F64 block so every row (including the ones that previously had
# no F cell at all) shares the same look the author already used for F5:F7
# etc.
$ws.Range("F5").Copy()
$ws.Range("F2:F64").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$values = @{
  2  = 0.489
  3  = 0.448
  4  = 0.469
  5  = 0.429
  6  = 0.43
  7  = 0.429
  8  = 0.477
  9  = 0.436
  10 = 0.48
  11 = 0.43
  12 = 0.443
  13 = 0.426
  14 = 0.442
  15 = 0.435
  16 = 0.381
  17 = 0.35
  18 = 0.38
  19 = 0.441
  20 = 0.191
  21 = 0.264
  22 = 0.2
  23 = 0.202
  24 = 0.212
  25 = 0.2
  26 = 0.194
  27 = 0.189
  28 = 0.164
  29 = 0.44
  30 = 0.45
  31 = 0.46
  32 = 0.47
  33 = 0.43
  34 = 0.414
  35 = 0.412
  36 = 0.39
  37 = 0.4
  38 = 0.425
  39 = 0.429
  40 = 0.423
  41 = 0.203
  42 = 0.202
  43 = 0.2
  44 = 0.203
  45 = 0.242
  46 = 0.151
  47 = 0.416
  48 = 0.445
  49 = 0.38
  50 = 0.44
  51 = 0.48
  52 = 0.424
  53 = 0.22
  54 = 0.2
  55 = 0.23
  56 = 0.18
  57 = 0.19
  58 = 0.16
  59 = 0.465
  60 = 0.47
  61 = 0.43
  62 = 0.46
  63 = 0.45
  64 = 0.5
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value2 = $values[$row]
}

# Match the author's final selection (cell F1) when they saved the file.
$ws.Range("F1").Select()
